$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the remaining "x" marks in the traceability matrix (Matriz de Trazabilidad CU vs Clases)
$ws.Range("G13").Value = "x"
$ws.Range("H13").Value = "x"

$ws.Range("E14").Value = "x"
$ws.Range("G14").Value = "x"
$ws.Range("H14").Value = "x"

$ws.Range("G15").Value = "x"
$ws.Range("H15").Value = "x"

$ws.Range("E16").Value = "x"
$ws.Range("G16").Value = "x"
$ws.Range("H16").Value = "x"

$ws.Range("G17").Value = "x"
$ws.Range("H17").Value = "x"

$ws.Range("G18").Value = "x"
$ws.Range("H18").Value = "x"

$ws.Range("E19").Value = "x"
$ws.Range("G19").Value = "x"
$ws.Range("H19").Value = "x"

# Update the active selection to match the saved workbook state
$ws.Range("K12").Select()
